$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of timesheet data (row 65 and 66), continuing the
# running-total formula in column C.

# Copy the date-column formatting down from row 64 so the new rows match
# the rest of the table (rather than creating a brand-new number format).
$ws.Range("A64").Copy()
$ws.Range("A65:A66").PasteSpecial(-4122)

$ws.Range("A65").Value = 45433
$ws.Range("B65").Value = 1.5

$ws.Range("A66").Value = 45434
$ws.Range("B66").Value = 0.5

# Fill the running-total formula into the two new rows as one operation so
# they form their own shared-formula group, same as Excel does when you
# drag/fill a formula down across several cells at once.
$ws.Range("C65:C66").Formula = "=C64+B65"

# Move selection to the new last cell, matching the updated navmenu target.
$ws.Range("C66").Select()
